$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.875.25'
$ws.Range('E2').Value = '  -4.63%  '
$ws.Range('D3').Value = '2.487.60'
$ws.Range('E3').Value = '  -3.36%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Formula = '="533.69"'
$ws.Range('D5').Copy() | Out-Null
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  -2.47%  '
$ws.Range('D6').Formula = '="142.55"'
$ws.Range('D6').Copy() | Out-Null
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('E6').Value = '  -7.55%  '
$ws.Range('D7').Formula = '="0.996"'
$ws.Range('D7').Copy() | Out-Null
$ws.Range('D7').PasteSpecial(-4163)
$ws.Range('E7').Value = '  -0.33%  '
$ws.Range('E8').Value = '  -4.11%  '
$ws.Range('D9').Value = '2.519.24'
$ws.Range('E9').Value = '  -2.19%  '
$ws.Range('E10').Value = '  -4.16%  '
$ws.Range('E11').Value = '  -2.65%  '
$ws.Range('D12').Formula = '="5.48"'
$ws.Range('D12').Copy() | Out-Null
$ws.Range('D12').PasteSpecial(-4163)
$ws.Range('E12').Value = '  +1.81%  '
$ws.Range('E13').Value = '  -3.66%  '
$ws.Range('D14').Value = '2.933.53'
$ws.Range('E14').Value = '  -3.11%  '
$ws.Range('E15').Value = '  -6.57%  '
$ws.Range('D16').Value = '58.753.10'
$ws.Range('E16').Value = '  -4.68%  '
$ws.Range('E17').Value = '  -4.00%  '
$ws.Range('D18').Value = '2.516.43'
$ws.Range('E18').Value = '  -2.35%  '
$ws.Range('D19').Formula = '="11.38"'
$ws.Range('D19').Copy() | Out-Null
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Value = '  -1.04%  '
$ws.Range('D20').Formula = '="4.27"'
$ws.Range('D20').Copy() | Out-Null
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('E20').Value = '  -5.81%  '
$ws.Range('D21').Formula = '="321.90"'
$ws.Range('D21').Copy() | Out-Null
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('E21').Value = '  -4.34%  '
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('E23').Value = '  -5.15%  '
$ws.Range('D24').Formula = '="60.87"'
$ws.Range('D24').Copy() | Out-Null
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Value = '  -3.67%  '
$ws.Range('D25').Formula = '="0.438"'
$ws.Range('D25').Copy() | Out-Null
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  -10.98%  '
$ws.Range('D26').Value = '2.625.07'
$ws.Range('E26').Value = '  -2.57%  '
$ws.Range('D27').Formula = '="0.997"'
$ws.Range('D27').Copy() | Out-Null
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('E27').Value = '  -0.31%  '
$ws.Range('E28').Value = '  -3.70%  '
$ws.Range('D29').Formula = '="7.70"'
$ws.Range('D29').Copy() | Out-Null
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Value = '  -4.46%  '
$ws.Range('D30').Formula = '="6.83"'
$ws.Range('D30').Copy() | Out-Null
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('E30').Value = '  -8.70%  '
$ws.Range('D31').Value = '0.0₃0762'
$ws.Range('E31').Value = '  -9.02%  '
$ws.Range('E32').Value = '  -3.70%  '
$ws.Range('E33').Value = '  -6.07%  '
$ws.Range('D34').Formula = '="0.996"'
$ws.Range('D34').Copy() | Out-Null
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('E34').Value = '  -0.31%  '
$ws.Range('D35').Formula = '="157.11"'
$ws.Range('D35').Copy() | Out-Null
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('E35').Value = '  -2.10%  '
$ws.Range('E36').Value = '  +0.64%  '
$ws.Range('D37').Formula = '="18.50"'
$ws.Range('D37').Copy() | Out-Null
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('E37').Value = '  -3.31%  '
$ws.Range('D38').Formula = '="4.37"'
$ws.Range('D38').Copy() | Out-Null
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('E38').Value = '  -7.54%  '
$ws.Range('E39').Value = '  -10.11%  '
$ws.Range('E40').Value = '  +1.20%  '
$ws.Range('D41').Formula = '="309.23"'
$ws.Range('D41').Copy() | Out-Null
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('E41').Value = '  -6.89%  '
$ws.Range('D42').Formula = '="36.79"'
$ws.Range('D42').Copy() | Out-Null
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('E42').Value = '  -1.93%  '
$ws.Range('E43').Value = '  -7.35%  '
$ws.Range('D44').Formula = '="0.791"'
$ws.Range('D44').Copy() | Out-Null
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  -13.78%  '
$ws.Range('D45').Formula = '="0.995"'
$ws.Range('D45').Copy() | Out-Null
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Value = '  -0.27%  '
$ws.Range('D46').Formula = '="0.595"'
$ws.Range('D46').Copy() | Out-Null
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  -1.60%  '
$ws.Range('D47').Formula = '="10.78"'
$ws.Range('D47').Copy() | Out-Null
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('E47').Value = '  -1.38%  '
$ws.Range('D48').Formula = '="124.29"'
$ws.Range('D48').Copy() | Out-Null
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Value = '  +1.21%  '
$ws.Range('E49').Value = '  -4.27%  '
$ws.Range('D50').Formula = '="18.60"'
$ws.Range('D50').Copy() | Out-Null
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('E50').Value = '  -4.81%  '
$ws.Range('E51').Value = '  -5.56%  '

$excel.CutCopyMode = 0

